$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsNote = $wb.Worksheets.Item("note")

# Update the category label for rows 58-120 on the "data" sheet:
# previously mislabeled as "l. Baumwoll - Spinnereien im Jahre 1854." (cotton spinning),
# correct it to a new category "m. Papier - Erzeugung im Jahre 1854." (paper production).
$wsData.Range("A58:A120").Value = "m. Papier - Erzeugung im Jahre 1854."

# Update view/selection state: "data" becomes the active/selected sheet,
# and its frozen-pane view scrolls down near the bottom of the data,
# while "note" is no longer the tab-selected sheet.
$wsData.Activate()

$wsData.Range("B153").Select()
$excel.ActiveWindow.ScrollRow = 119

$wsNote.Range("A2").Select()
